# Apply coin-ranking snapshot refresh: updated prices/volumes and row-shifted coin list (rows 7-17).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + "300.80"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'" + "-3.08%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'" + "35.48"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'" + "0.20%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'" + "5.068"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'" + "-0.71%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'" + "0.07982"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'" + "-2.67%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'" + "1.896"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'" + "-8.06%"
$ws.Range("E6").Style = "Normal"
$ws.Range("B7").Value = "KuCoinToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D7").Value = "'" + "7.755"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'" + "-2.40%"
$ws.Range("E7").Style = "Normal"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "'" + "0.9283"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'" + "0.26%"
$ws.Range("E8").Style = "Normal"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").Value = "'" + "0.1455"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'" + "31.41%"
$ws.Range("E9").Style = "Normal"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'" + "0.1904"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'" + "-0.49%"
$ws.Range("E10").Style = "Normal"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'" + "0.09032"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'" + "-2.24%"
$ws.Range("E11").Style = "Normal"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'" + "0.03491"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'" + "-4.20%"
$ws.Range("E12").Style = "Normal"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'" + "0.09856"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'" + "-0.54%"
$ws.Range("E13").Style = "Normal"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'" + "0.001395"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'" + "-2.23%"
$ws.Range("E14").Style = "Normal"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "'" + "0.005721"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'" + "-2.13%"
$ws.Range("E15").Style = "Normal"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "'" + "3.531"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'" + "1.54%"
$ws.Range("E16").Style = "Normal"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").Value = "'" + "4.055"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'" + "-1.68%"
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'" + "2.89%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'" + "0.3447"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'" + "1.15%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'" + "0.1304"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'" + "-0.33%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'" + "5.032"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'" + "-1.30%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'" + "0.2399"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'" + "8.91%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'" + "0.04497"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'" + "-1.07%"
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'" + "-1.01%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'" + "0.004758"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'" + "-1.17%"
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'" + "-1.59%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'" + "0.0003027"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'" + "-31.87%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D39").Value = "'" + "0.01822"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'" + "-7.67%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'" + "0.04738"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'" + "-2.81%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'" + "0.01057"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'" + "16.65%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'" + "0.007324"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'" + "-3.79%"
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'" + "-4.25%"
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'" + "-3.22%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'" + "0.01094"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'" + "-5.99%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'" + "0.00006221"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'" + "-4.88%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'" + "0.00000000751"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'" + "0.21%"
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'" + "6.10%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D50").Value = "'" + "0.00002102"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'" + "0.21%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'" + "0.0002002"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'" + "0.21%"
$ws.Range("E51").Style = "Normal"
